$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.220.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").Value = "'2.243.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.14%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'246.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("E6").Value = "  -3.40%  "
$ws.Range("D7").Value = "'74.84"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.75%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "'0.618"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.53%  "
$ws.Range("D10").Value = "'42.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.37%  "
$ws.Range("E11").Value = "  -2.50%  "
$ws.Range("D12").Value = "'7.18"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.29%  "
$ws.Range("E13").Value = "  -1.52%  "
$ws.Range("D14").Value = "'14.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.41%  "
$ws.Range("D15").Value = "'0.853"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.46%  "
$ws.Range("D16").Value = "'2.272.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("D17").Value = "'42.099.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.24%  "
$ws.Range("D18").Value = "'0.0₃0982"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").Value = "'6.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("E21").Value = "  +4.62%  "
$ws.Range("D22").Value = "'231.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("E23").Value = "  +42.16%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "'11.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.17%  "
$ws.Range("D26").Value = "'3.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.01%  "
$ws.Range("E27").Value = "  -2.56%  "
$ws.Range("D28").Value = "'2.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("D29").Value = "'169.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.14%  "
$ws.Range("D30").Value = "'20.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("D31").Value = "'0.0822"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.56%  "
$ws.Range("D32").Value = "'31.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.72%  "
$ws.Range("E33").Value = "  -1.79%  "
$ws.Range("D34").Value = "'0.124"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.95%  "
$ws.Range("D35").Value = "'5.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.98%  "
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("E37").Value = "  +2.45%  "
$ws.Range("D38").Value = "'13.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("D39").Value = "'2.18"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.93%  "
$ws.Range("E40").Value = "  -1.11%  "
$ws.Range("D41").Value = "'62.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.42%  "
$ws.Range("D42").Value = "'0.205"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.28%  "
$ws.Range("D43").Value = "'106.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.96%  "
$ws.Range("E44").Value = "  +2.34%  "
$ws.Range("D45").Value = "'8.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("D46").Value = "'0.997"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("E47").Value = "  -2.45%  "
$ws.Range("D48").Value = "'1.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("E49").Value = "  -4.84%  "
$ws.Range("D50").Value = "'2.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.93%  "
$ws.Range("D51").Value = "'4.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.44%  "
